$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 19.02.2022 15:45"

# Row 3 (Tesco): new price scraped -> shift current price to "old", store new price,
# and write the delta/old-date columns as plain text (as produced by the scraping script)
$ws.Range("C3").Value = 36.5
$ws.Range("B3").Value = 36.7

# D3 must hold the literal text "+0.2" (not be re-interpreted as a number) while keeping
# its original default (unstyled) cell format.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+0.2"
$ws.Range("D3").Style = "Normal"

# E3 must hold the literal text timestamp and lose its previous date number format.
$ws.Range("E3").Value = "2022-02-19 15:46:25"
$ws.Range("E3").Style = "Normal"
